$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 44320
$ws.Range("M3").Value = 50
$ws.Range("N3").Value = 18000
$ws.Range("P3").Value = 18800
$ws.Range("R3").Value = "Provincia de Limarí"
$ws.Range("S3").Value = 1044

$ws.Range("D4").Value = 44362
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 19000
$ws.Range("P4").Value = 19500
$ws.Range("R4").Value = "Provincia de Curicó"
$ws.Range("S4").Value = 1083
